# Quad Report weekly refresh: bump the "Week Ending" date shown in the
# report header from 2016-07-25 to 2016-08-01.
#
# The header text lives in a single paragraph built from several runs;
# the date is the tail of the last run ("...Week Ending:  2016-07-25").
# We locate that substring and overwrite just the date portion through
# TextRange.Characters(start, length), which is how PowerPoint itself
# splits a run when only part of its text is edited interactively -
# the newly typed text becomes its own run, inheriting the same
# character formatting (rPr) as the run it was carved out of.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$oldDate = "2016-07-25"
$newDate = "2016-08-01"

$found = $false
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    $tf = $shp.TextFrame
    if (-not $tf.HasText) { continue }
    $tr = $tf.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($oldDate)
    if ($idx -ge 0) {
        $sub = $tr.Characters($idx + 1, $oldDate.Length)
        $sub.Text = $newDate
        $found = $true
    }
}

if (-not $found) {
    throw "Could not find '$oldDate' on slide 1 to update."
}
